$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.816264271736145
$ws.Range("B1").Value = 1.865287899971008
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.97567343711853
$ws.Range("E1").Value = 0.6954168677330017
